# "10 years Finalization data" -- duplicate the finalized daily-data table
# (the "Data Harian - Table" sheet's A9:K40 block, i.e. header + 31 days of
# July 2022 readings) onto a brand-new sheet, mirroring Excel's
# "Move or Copy Sheet... (Create a copy)" / copy-paste workflow that the
# author used to produce a clean, standalone finalized-data sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# New sheet is inserted right after the existing one, becoming "Sheet1".
$newSheet = $wb.Worksheets.Add($null, $ws1)

# Copy the table (header row + 31 daily rows) including all formatting
# (borders, alignment, wrap) down to the top-left of the new sheet.
$srcRange  = $ws1.Range("A9:K40")
$destRange = $newSheet.Range("A1")
$srcRange.Copy($destRange)

# Leave the original sheet's view focused on the table that was copied.
$ws1.Activate() | Out-Null
$ws1.Range("A9:K40").Select() | Out-Null
$excel.ActiveWindow.DisplayGridlines = $true

# Finish with the new sheet active/selected, matching the tab the author
# ended up on after the copy.
$newSheet.Activate() | Out-Null
$newSheet.Range("A1:K32").Select() | Out-Null
